# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# Updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 12;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 16;  Tag = "ba"; Act = "Appreciation" },
    @{ Row = 28;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 32;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 34;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 45;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 46;  Tag = "sv"; Act = "Statement-opinion" },
    @{ Row = 50;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 52;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 60;  Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 62;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 66;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 68;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 89;  Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 92;  Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 96;  Tag = "b";  Act = "Acknowledge (Backchannel)" },
    @{ Row = 100; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 101; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 111; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 114; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 120; Tag = "ba"; Act = "Appreciation" },
    @{ Row = 133; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 155; Tag = "sd"; Act = "Statement-non-opinion" },
    @{ Row = 168; Tag = "%";  Act = "Uninterpretable" },
    @{ Row = 172; Tag = "aa"; Act = "Agree/Accept" },
    @{ Row = 187; Tag = "qy"; Act = "Yes-No-Question" }
)

foreach ($u in $updates) {
    $ws.Range("I" + $u.Row).Value = $u.Tag
    $ws.Range("J" + $u.Row).Value = $u.Act
}
